$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (BOM gained a "SUP75P03-07-E3" P-Channel MOSFET line).
# This shifts the old rows 17-22 down to 18-23 and carries their cell
# formatting/styles along automatically (same as Excel's real Insert behavior).
$ws.Rows(17).Insert()

# --- Row 2: Polarized Capacitor (1000uF) - reference designator renumbered
$ws.Range("E2").Value = "C11"

# --- Row 3: Ceramic Capacitor (10uF) - reference designators renumbered
$ws.Range("E3").Value = "C8, C9"

# --- Row 4: Ceramic Capacitor (1uF) - qty/refs renumbered, note updated
$ws.Range("A4").Value = 3
$ws.Range("E4").Value = "C4, C5, C10"
$ws.Range("H4").Value = "Bulk (maybe)"

# --- Row 5: Ceramic Capacitor (0.1uF) - qty/refs renumbered
$ws.Range("A5").Value = 6
$ws.Range("E5").Value = "C1, C2, C7"

# --- Row 6: Ceramic Capacitor (20pF) - refs renumbered
$ws.Range("E6").Value = "C3, C6"

# --- Row 7: Resistor - value/refs/part number changed (330K -> 560K)
$ws.Range("B7").Value = "560K"
$ws.Range("E7").Value = "R7"
$ws.Range("F7").Value = "MFR-25FBF52-576K"

# --- Row 8: Resistor (10K) - refs renumbered
$ws.Range("E8").Value = "R2, R8"

# --- Row 9: Resistor (1K) - refs renumbered
$ws.Range("E9").Value = "R5, R6"

# --- Row 10: Resistor (220) - refs renumbered
$ws.Range("E10").Value = "R3, R4"

# --- Row 11: Resistor (68) - refs renumbered
$ws.Range("E11").Value = "R1"

# --- Row 12: Zener Diode (27V) - refs swapped
$ws.Range("E12").Value = "D2"

# --- Row 13: Zener Diode (22V) - refs swapped
$ws.Range("E13").Value = "D1"

# --- Row 14: Zener Diode (10V) - refs renumbered
$ws.Range("E14").Value = "D3, D4"

# --- Row 16: IPP80P03P4L-04 P-Channel MOSFET - qty corrected
$ws.Range("A16").Value = 1

# --- Row 17 (new): SUP75P03-07-E3 P-Channel MOSFET
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "SUP75P03-07-E3"
$ws.Range("C17").Value = "SUP75P03-07-E3-H"
$ws.Range("D17").Value = "TO220BH"
$ws.Range("E17").Value = "Q2"
$ws.Range("F17").Value = "SUP75P03-07-E3"
$ws.Range("G17").Value = "P-Channel MOSFET"

# --- Row 21 (was 20): Crystal Oscillator - ref designator added
$ws.Range("E21").Value = "Y1"

# --- Row 22 (was 21): Fuse - note spelling fix ("experiement" -> "experiment")
$ws.Range("H22").Value = "Need to experiment to find appropriate fuse type"

# Restore the author's active-cell selection at save time
$ws.Range("D8").Select() | Out-Null
